$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("airline_financials")

# Row 194: 2025 Q4 DAL
$ws.Range("A194").Value = 2025
$ws.Range("B194").Value = 4
$ws.Range("C194").Value = "DAL"
$ws.Range("D194").Value = 16003000000
$ws.Range("E194").Value = 12916000000
$ws.Range("F194").Value = 14536000000
$ws.Range("G194").Value = 1219000000
$ws.Range("H194").Value = 59861000000
$ws.Range("I194").Value = 72946000000
$ws.Range("K194").Value = 351000000

# Row 195: 2025 FY DAL
$ws.Range("A195").Value = 2025
$ws.Range("B195").Value = "FY"
$ws.Range("C195").Value = "DAL"
$ws.Range("D195").Value = 63364000000
$ws.Range("E195").Value = 51768000000
$ws.Range("F195").Value = 57542000000
$ws.Range("G195").Value = 5005000000
$ws.Range("H195").Value = 249578000000
$ws.Range("I195").Value = 298045000000
$ws.Range("K195").Value = 1337000000

# Row 196: 2025 Q4 UAL (pending release)
$ws.Range("A196").Value = 2025
$ws.Range("B196").Value = 4
$ws.Range("C196").Value = "UAL"

# Row 197: 2025 FY UAL (pending release)
$ws.Range("A197").Value = 2025
$ws.Range("B197").Value = "FY"
$ws.Range("C197").Value = "UAL"

# Row 198: 2025 Q4 LUV (pending release)
$ws.Range("A198").Value = 2025
$ws.Range("B198").Value = 4
$ws.Range("C198").Value = "LUV"

# Row 199: 2025 FY LUV (pending release)
$ws.Range("A199").Value = 2025
$ws.Range("B199").Value = "FY"
$ws.Range("C199").Value = "LUV"

# Row 200: 2025 Q4 AAL (pending release)
$ws.Range("A200").Value = 2025
$ws.Range("B200").Value = 4
$ws.Range("C200").Value = "AAL"

# Row 201: 2025 FY AAL (pending release)
$ws.Range("A201").Value = 2025
$ws.Range("B201").Value = "FY"
$ws.Range("C201").Value = "AAL"

# Update view to reflect scroll/selection near the new rows
$ws.Range("G194").Select()
